$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B (Tempo (ms)) updates for rows 2-11
$ws.Range("B2").Value = 63
$ws.Range("B3").Value = 64
$ws.Range("B4").Value = 64
$ws.Range("B5").Value = 64
$ws.Range("B6").Value = 65
$ws.Range("B7").Value = 66
$ws.Range("B8").Value = 68
$ws.Range("B9").Value = 69
$ws.Range("B10").Value = 69
$ws.Range("B11").Value = 70

# Column C (Memória (KB)) update for row 11
$ws.Range("C11").Value = 3300.171875

# Summary rows (stored as text strings)
$ws.Range("B13").NumberFormat = "@"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("B14").NumberFormat = "@"

$ws.Range("B13").Value = "66.20"
$ws.Range("C13").Value = "-827.48"
$ws.Range("B14").Value = "65.50"
